# Insert a new row at position 144 (shifting existing rows 144-234 down to 145-235)
# and populate the new row 144 with the new data record ("Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 144, shifting rows 144:234 down to 145:235
$ws.Rows.Item(144).Insert()

# Copy formats (and values) from row 145 (which now holds what used to be row 144)
# down into the freshly inserted, blank row 144 so formatting - especially the
# date number-format on column D - matches the rest of the table.
$ws.Range("A145:R145").Copy()
$ws.Range("A144:R144").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Now overwrite the new row 144 with its actual values.
$ws.Cells.Item(144, 1).Value = 5
$ws.Cells.Item(144, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(144, 3).Value = "Maule"
$ws.Cells.Item(144, 4).Value = 44582
$ws.Cells.Item(144, 5).Value = 7
$ws.Cells.Item(144, 6).Value = 100112003
$ws.Cells.Item(144, 7).Value = "Ajo"
$ws.Cells.Item(144, 8).Value = "Chino"
$ws.Cells.Item(144, 9).Value = "Primera"
$ws.Cells.Item(144, 10).Value = 150
$ws.Cells.Item(144, 11).Value = 20000
$ws.Cells.Item(144, 12).Value = 20000
$ws.Cells.Item(144, 13).Value = 20000
$ws.Cells.Item(144, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(144, 15).Value = "China"
$ws.Cells.Item(144, 16).Value = 2000
$ws.Cells.Item(144, 17).Value = 10
$ws.Cells.Item(144, 18).Value = "Hortaliza"
